# feat: add 2022-Q1 data
#
# The workbook's last sheet "总计" (rolling summary) is duplicated first so
# the duplicate can keep serving as the "总计" sheet (with a new leading
# row for 2022-Q1 inserted), while the original "总计" sheet is renamed to
# "2022-Q1" and repurposed to hold the new quarter's per-fund holdings
# detail (same layout as the other quarter sheets).

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1. Duplicate the summary sheet now (while it still has the old data +
#    formatting) so the copy can become the refreshed "总计" sheet at the
#    end of the workbook, after everything else.
# ---------------------------------------------------------------------
$summary.Copy($null, $summary) | Out-Null
$total = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 2. Repurpose the original sheet -> "2022-Q1" quarter-detail sheet
#    (rename the original out of the way first so the duplicate is free
#    to take over the "总计" name).
# ---------------------------------------------------------------------
$q1 = $summary
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

$total.Name = "总计"

# Use an existing quarter sheet as the formatting template (header row +
# index column both use the bold/bordered/centered style already defined
# in the workbook).
$template = $wb.Worksheets.Item("2021-Q4")

$template.Range("B1:H1").Copy() | Out-Null
$q1.Range("B1:H1").PasteSpecial(-4122) | Out-Null

$template.Range("A2").Copy() | Out-Null
$q1.Range("A2:A18").PasteSpecial(-4122) | Out-Null

# Headers
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Fund code / size / position columns are numeric-looking strings (leading
# zeros in fund codes, fixed 2dp strings for the figures) that must be
# forced to text storage so they aren't auto-coerced to numbers, matching
# the other quarter sheets. The fund-name column (C) is already non-numeric
# text and needs no special handling.
$q1.Range("B2:B18").NumberFormat = "@"
$q1.Range("D2:G18").NumberFormat = "@"

$rows = @(
    @(0,  "000362", "国泰聚信价值优势灵活配置混合A", "56.15", "89.00", "5.02", "2.8187", 5),
    @(1,  "008415", "国泰大制造两年持有期混合",      "23.19", "92.05", "5.35", "1.2407", 4),
    @(2,  "012173", "国泰兴泽优选一年持有期混合A",    "16.89", "89.83", "5.21", "0.8800", 4),
    @(3,  "000363", "国泰聚信价值优势灵活配置混合C", "17.09", "89.00", "5.02", "0.8579", 5),
    @(4,  "020010", "国泰金牛创新混合",              "16.99", "84.21", "4.51", "0.7662", 5),
    @(5,  "007835", "国泰鑫睿混合",                  "9.37",  "78.94", "4.77", "0.4469", 3),
    @(6,  "012174", "国泰兴泽优选一年持有期混合C",    "7.14",  "89.83", "5.21", "0.3720", 4),
    @(7,  "001579", "国泰大农业股票",                "11.72", "90.32", "3.05", "0.3575", 10),
    @(8,  "005244", "国泰聚优价值灵活配置混合A",      "7.72",  "83.97", "3.54", "0.2733", 4),
    @(9,  "005245", "国泰聚优价值灵活配置混合C",      "4.52",  "83.97", "3.54", "0.1600", 4),
    @(10, "003516", "国泰融安多策略灵活配置混合",      "11.18", "71.30", "1.43", "0.1599", 9),
    @(11, "001922", "国泰多策略收益灵活配置混合",      "6.86",  "24.58", "1.24", "0.0851", 1),
    @(12, "002197", "国泰鑫策略价值灵活配置混合",      "6.92",  "21.00", "1.05", "0.0727", 1),
    @(13, "001850", "国泰安益灵活配置混合A",          "6.63",  "21.23", "0.99", "0.0656", 1),
    @(14, "000367", "国泰安康定期支付混合A",          "5.28",  "21.92", "1.16", "0.0612", 1),
    @(15, "002061", "国泰安康定期支付混合C",          "2.35",  "21.92", "1.16", "0.0273", 1),
    @(16, "004252", "国泰安益灵活配置混合C",          "2.10",  "21.23", "0.99", "0.0208", 1)
)

$r = 2
foreach ($row in $rows) {
    $q1.Cells.Item($r, 1).Value = $row[0]
    $q1.Cells.Item($r, 2).Value = $row[1]
    $q1.Cells.Item($r, 3).Value = $row[2]
    $q1.Cells.Item($r, 4).Value = $row[3]
    $q1.Cells.Item($r, 5).Value = $row[4]
    $q1.Cells.Item($r, 6).Value = $row[5]
    $q1.Cells.Item($r, 7).Value = $row[6]
    $q1.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3. Refresh the (duplicated) "总计" sheet: insert a new leading row for
#    2022-Q1 and shift the old rows (and their A-column index) down one.
# ---------------------------------------------------------------------
$total.Range("A2").EntireRow.Insert(-4121)
$total.Range("A2:D2").ClearFormats()

$template.Range("A2").Copy() | Out-Null
$total.Range("A2").PasteSpecial(-4122) | Out-Null

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 17
$total.Range("D2").Value = 8.67

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
